$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 174, pushing the existing Rabanito (Vega Central
# Mapocho de Santiago) entries that previously lived in rows 174-183
# down to rows 175-184.
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with the new weekly price record.
$ws.Cells.Item(174, 1).Value = 9
$ws.Cells.Item(174, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(174, 3).Value = "Metropolitana"
$ws.Cells.Item(174, 4).Value = 44516
$ws.Cells.Item(174, 5).Value = 13
$ws.Cells.Item(174, 6).Value = 300000001
$ws.Cells.Item(174, 7).Value = "Rabanito"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 7900
$ws.Cells.Item(174, 11).Value = 2500
$ws.Cells.Item(174, 12).Value = 3000
$ws.Cells.Item(174, 13).Value = 2747
$ws.Cells.Item(174, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(174, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(174, 16).Value = 27
$ws.Cells.Item(174, 17).Value = 100
$ws.Cells.Item(174, 18).Value = "Hortaliza"
